# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Sat Jan 27 11:33:27 UTC 2024 with GitHub Actions".
#
# Column D ("Price") cells hold numeric-looking text (e.g. "91.90", "0.999")
# that must stay as literal text, matching digit-for-digit what a scraper
# would have written (trailing zeros, multi-dot thousands separators, etc.).
# Assigning a numeric-looking string straight to Range.Value lets Excel
# auto-coerce it to a real number (dropping trailing zeros / losing exact
# text), so we force the cell to Text format first for every Price update.
# Column E/B/C values are never ambiguous (contain "%", letters or URLs) so
# a plain assignment is safe and keeps them as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.727.05"
$ws.Range("E2").Value = "  +1.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.69"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.03"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.90"
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("E7").Value = "  +1.87%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  -0.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.29"
$ws.Range("E10").Value = "  +0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.38"
$ws.Range("E11").Value = "  +1.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0796"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("E13").Value = "  -0.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.63"
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.614.29"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16
$ws.Range("E16").Value = "  +0.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.251.04"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("E18").Value = "  +1.69%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.640.25"
$ws.Range("E19").Value = "  +1.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  +6.27%  "

# Row 21
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").Value = "  +1.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.92"
$ws.Range("E23").Value = "  +0.50%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.72"
$ws.Range("E24").Value = "  -0.13%  "

# Row 25
$ws.Range("E25").Value = "  +1.24%  "

# Row 27
$ws.Range("E27").Value = "  +2.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.97"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("E29").Value = "  +0.23%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.45"
$ws.Range("E31").Value = "  +2.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.90"
$ws.Range("E32").Value = "  +1.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.26"
$ws.Range("E33").Value = "  +3.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("E35").Value = "  +1.93%  "

# Row 36
$ws.Range("E36").Value = "  -1.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.83"
$ws.Range("E37").Value = "  +2.19%  "

# Row 38
$ws.Range("E38").Value = "  +0.95%  "

# Row 39
$ws.Range("E39").Value = "  +1.31%  "

# Row 40
$ws.Range("E40").Value = "  -0.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  +0.17%  "

# Row 42
$ws.Range("E42").Value = "  +0.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.025.23"
$ws.Range("E43").Value = "  -3.30%  "

# Row 44
$ws.Range("E44").Value = "  -2.91%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.43"
$ws.Range("E45").Value = "  +1.60%  "

# Row 46
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("E47").Value = "  +14.13%  "

# Row 48
$ws.Range("E48").Value = "  -1.45%  "

# Row 49
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.69"
$ws.Range("E49").Value = "  +4.11%  "

# Row 50
$ws.Range("E50").Value = "  +1.21%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("E51").Value = "  -1.40%  "
